# Apply pylinting/readme-related data fix:
# - Row 25 cells get the "normal" style (s="6") like the rest of the table rows.
# - A new row 26 is appended, duplicating the "Sem" skill-score row (same
#   values as row 3), left unstyled like row 25 was before this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output")

# --- Fix up row 25 styling: give every populated cell (A25:U25) the same
#     style used throughout the rest of the table (style index 6, i.e. the
#     sheet's default/general style).
$styledRow = $ws.Range("A24:U24")
$targetRow = $ws.Range("A25:U25")
$targetRow.Style = $styledRow.Style

# --- Append row 26 with the same data set found in row 3 ("Sem" skill
#     scores), which has no explicit style (matches the pre-edit row 25).
$row26Values = @(
    "99004351", "Sem", 93, 52, 65, 39, 87, 43, 55, 65, 99, 91, 64, 75, 79, 86, 44, 75, 98, 23, 91
)

$col = 1
foreach ($val in $row26Values) {
    $ws.Cells.Item(26, $col).Value = $val
    $col++
}
